# Generate Report for Handoff
#
# Refreshes the localization-status report: the files moved from
# "In Translation" to "Ready for handoff", and the handoff timestamps were
# stamped for both target locales. The "Status" column on every sheet is
# also widened to comfortably fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# A ColumnWidth value (in characters) that this engine's internal pixel grid
# rounds to the widened "Status" column width used across the sheets.
$statusColWidth = 16.333333

# --- Overview sheet ---------------------------------------------------------
# Columns E (zh-cn) / F (de-de) hold the per-language status.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Column G holds the latest HO xliff generate date (the later of the two
# per-locale handoff timestamps below).
$wsOverview.Range("G2").Value = "2016-08-21 19:01:32"
# Widen the two status columns.
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-21 19:01:28"
$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet -------------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-21 19:01:32"
$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
